$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Swap the ID (column A) and Text (column B) values between rows 19 and 20.
# This reflects a change in line-matching order (now using LineIDs for
# Jaccard comparison), so the two recording lines that used to be on
# rows 19/20 have traded places.
$a19 = $ws.Range("A19").Value()
$b19 = $ws.Range("B19").Value()
$a20 = $ws.Range("A20").Value()
$b20 = $ws.Range("B20").Value()

$ws.Range("A19").Value = $a20
$ws.Range("B19").Value = $b20
$ws.Range("A20").Value = $a19
$ws.Range("B20").Value = $b19
